# Commit: "Learning by not doing analysis & replace eff_cost_loan with APR"
#
# The underlying regression results (produced upstream, outside this
# workbook) were re-run and a handful of coefficients/standard errors in
# the "First treatment" columns (source columns D/E -> displayed columns
# E/F) changed. Update the displayed table's cached numbers to match the
# refreshed analysis.
#
# These cells hold formulas that pull their cached text from an external
# workbook link ([1]multiple_loans!...). That external source file isn't
# reachable from this environment (no real linked workbook to refresh
# against), so instead of leaving the figures stale we write the new,
# already-known figures directly onto the cells - using a literal-text
# formula (="...") rather than a plain value so the cell keeps its
# "str" formula-result typing/style rather than turning into a bare
# number/shared-string, matching how these text-formatted regression
# figures (with parens/stars) are stored elsewhere in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Formula  = '="-373.8**"'
$ws.Range("F5").Formula  = '="-50.3***"'

$ws.Range("E6").Formula  = '="(148.3)"'

$ws.Range("E7").Formula  = '="-105.1"'
$ws.Range("F7").Formula  = '="0.19"'

$ws.Range("E8").Formula  = '="(146.4)"'
$ws.Range("F8").Formula  = '="(7.79)"'

$ws.Range("E11").Formula = '="0.007"'
$ws.Range("F11").Formula = '="0.039"'
